# The edit rotates the full record content of two groups of three rows on
# the "Artfynd" sheet:
#   - rows 10, 11, 12  ->  row 10 gets what row 11 had, row 11 gets what row
#     12 had, and row 12 gets what row 10 had (a 3-row cyclic shift).
#   - rows 28, 29, 30  ->  same rotation pattern (29->28, 30->29, 28->30).
#
# Only the cells that actually differ between the old and new record for each
# row are written below (columns that stay identical, e.g. S/T/U/V/W/Y/AA/
# AD/AE/AG/AT/AW/AX/AY, are left untouched).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Row 10  (becomes the old row 11 record: Garnlav / Alectoria sarmentosa)
# ---------------------------------------------------------------------
$ws.Range("A10").Value = 131244279
$ws.Range("B10").Value = 79244
$ws.Range("E10").Value = 6425
$ws.Range("F10").Value = "Garnlav"
$ws.Range("G10").Value = "Alectoria sarmentosa"
$ws.Range("H10").Value = "(Ach.) Ach."
$ws.Range("J10").Value = ""
$ws.Range("L10").ClearContents()
$ws.Range("M10").ClearContents()
$ws.Range("Q10").Value = 613427
$ws.Range("R10").Value = 6998062
$ws.Range("AF10").Value = ""
$ws.Range("AJ10").Value = "gran"
$ws.Range("AK10").Value = "Picea abies"
$ws.Range("AO10").Value = "Picea abies"

# ---------------------------------------------------------------------
# Row 11  (becomes the old row 12 record: Garnlav / Alectoria sarmentosa,
#          but on a "tall" (pine) substrate)
# ---------------------------------------------------------------------
$ws.Range("A11").Value = 131244300
$ws.Range("Q11").Value = 613444
$ws.Range("R11").Value = 6998046
$ws.Range("AJ11").Value = "tall"
$ws.Range("AK11").Value = "Pinus sylvestris"
$ws.Range("AO11").Value = "Pinus sylvestris"

# ---------------------------------------------------------------------
# Row 12  (becomes the old row 10 record: Spillkråka / Dryocopus martius)
# ---------------------------------------------------------------------
$ws.Range("A12").Value = 131244258
$ws.Range("B12").Value = 57881
$ws.Range("E12").Value = 100049
$ws.Range("F12").Value = "Spillkråka"
$ws.Range("G12").Value = "Dryocopus martius"
$ws.Range("H12").Value = "(Linnaeus, 1758)"
$ws.Range("J12").ClearContents()
$ws.Range("L12").Value = ""
$ws.Range("M12").Value = "färska spår"
$ws.Range("Q12").Value = 613451
$ws.Range("R12").Value = 6998138
$ws.Range("AF12").ClearContents()
$ws.Range("AJ12").ClearContents()
$ws.Range("AK12").ClearContents()
$ws.Range("AO12").ClearContents()

# ---------------------------------------------------------------------
# Row 28  (becomes the old row 29 record: Rosenticka / Fomitopsis rosea)
# ---------------------------------------------------------------------
$ws.Range("A28").Value = 131244259
$ws.Range("B28").Value = 92107
$ws.Range("E28").Value = 658
$ws.Range("F28").Value = "Rosenticka"
$ws.Range("G28").Value = "Fomitopsis rosea"
$ws.Range("H28").Value = "(Alb. & Schwein.:Fr.) P.Karst."
$ws.Range("I28").Value = ""
$ws.Range("K28").ClearContents()
$ws.Range("L28").ClearContents()
$ws.Range("M28").ClearContents()
$ws.Range("N28").ClearContents()
$ws.Range("Q28").Value = 613387
$ws.Range("R28").Value = 6998216
$ws.Range("Z28").ClearContents()
$ws.Range("AB28").ClearContents()

# ---------------------------------------------------------------------
# Row 29  (becomes the old row 30 record: Tretåig hackspett / Picoides
#          tridactylus, with a public comment about fresh ring pecking)
# ---------------------------------------------------------------------
$ws.Range("A29").Value = 131244252
$ws.Range("B29").Value = 57884
$ws.Range("E29").Value = 100109
$ws.Range("F29").Value = "Tretåig hackspett"
$ws.Range("G29").Value = "Picoides tridactylus"
$ws.Range("H29").Value = "(Linnaeus, 1758)"
$ws.Range("K29").Value = ""
$ws.Range("L29").Value = ""
$ws.Range("M29").Value = "färska spår"
$ws.Range("N29").Value = ""
$ws.Range("Q29").Value = 613433
$ws.Range("R29").Value = 6998019
$ws.Range("AC29").Value = "Färska ringhack på tall"

# ---------------------------------------------------------------------
# Row 30  (becomes the old row 28 record: Talltita / Poecile montanus)
# ---------------------------------------------------------------------
$ws.Range("A30").Value = 131244255
$ws.Range("B30").Value = 58043
$ws.Range("E30").Value = 103021
$ws.Range("F30").Value = "Talltita"
$ws.Range("G30").Value = "Poecile montanus"
$ws.Range("H30").Value = "(Conrad von Baldenstein, 1827)"
# "4" must stay textual (Antal is stored as text in this sheet) - a leading
# quote forces Excel to keep it as text instead of coercing it to a number.
$ws.Range("I30").Value = "'4"
$ws.Range("M30").Value = "födosökande"
$ws.Range("N30").Value = "observerad"
$ws.Range("Q30").Value = 613399
$ws.Range("R30").Value = 6998197
$ws.Range("Z30").Value = "09:50"
$ws.Range("AB30").Value = "10:00"
$ws.Range("AC30").ClearContents()
